$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Bold header font on Catalogos sheet (A1:B1)
$ws2.Range("A1:B1").Font.Bold = $true

# Column widths for Catalogos sheet
$ws2.Columns.Item(1).ColumnWidth = 28.8
$ws2.Columns.Item(2).ColumnWidth = 15.9

# Countries list (column A, rows 2-193)
$ws2.Cells.Item(2, 1).Value = "Afganistán"
$ws2.Cells.Item(3, 1).Value = "Albania"
$ws2.Cells.Item(4, 1).Value = "Alemania"
$ws2.Cells.Item(5, 1).Value = "Andorra"
$ws2.Cells.Item(6, 1).Value = "Angola"
$ws2.Cells.Item(7, 1).Value = "Antigua y Barbuda"
$ws2.Cells.Item(8, 1).Value = "Arabia Saudita"
$ws2.Cells.Item(9, 1).Value = "Argelia"
$ws2.Cells.Item(10, 1).Value = "Argentina"
$ws2.Cells.Item(11, 1).Value = "Armenia"
$ws2.Cells.Item(12, 1).Value = "Australia"
$ws2.Cells.Item(13, 1).Value = "Azerbaiyán"
$ws2.Cells.Item(14, 1).Value = "Bahamas"
$ws2.Cells.Item(15, 1).Value = "Bangladés"
$ws2.Cells.Item(16, 1).Value = "Barbados"
$ws2.Cells.Item(17, 1).Value = "Baréin"
$ws2.Cells.Item(18, 1).Value = "Bélgica"
$ws2.Cells.Item(19, 1).Value = "Belice"
$ws2.Cells.Item(20, 1).Value = "Benín"
$ws2.Cells.Item(21, 1).Value = "Bielorrusia"
$ws2.Cells.Item(22, 1).Value = "Birmania"
$ws2.Cells.Item(23, 1).Value = "Bolivia"
$ws2.Cells.Item(24, 1).Value = "Bosnia y Herzegovina"
$ws2.Cells.Item(25, 1).Value = "Botsuana"
$ws2.Cells.Item(26, 1).Value = "Brasil"
$ws2.Cells.Item(27, 1).Value = "Brunéi"
$ws2.Cells.Item(28, 1).Value = "Bulgaria"
$ws2.Cells.Item(29, 1).Value = "Burkina Faso"
$ws2.Cells.Item(30, 1).Value = "Burundi"
$ws2.Cells.Item(31, 1).Value = "Bután"
$ws2.Cells.Item(32, 1).Value = "Cabo Verde"
$ws2.Cells.Item(33, 1).Value = "Camboya"
$ws2.Cells.Item(34, 1).Value = "Camerún"
$ws2.Cells.Item(35, 1).Value = "Canadá"
$ws2.Cells.Item(36, 1).Value = "Catar"
$ws2.Cells.Item(37, 1).Value = "Chad"
$ws2.Cells.Item(38, 1).Value = "Chile"
$ws2.Cells.Item(39, 1).Value = "China"
$ws2.Cells.Item(40, 1).Value = "Chipre"
$ws2.Cells.Item(41, 1).Value = "Ciudad del Vaticano"
$ws2.Cells.Item(42, 1).Value = "Colombia"
$ws2.Cells.Item(43, 1).Value = "Comoras"
$ws2.Cells.Item(44, 1).Value = "Corea del Norte"
$ws2.Cells.Item(45, 1).Value = "Corea del Sur"
$ws2.Cells.Item(46, 1).Value = "Costa de Marfil"
$ws2.Cells.Item(47, 1).Value = "Costa Rica"
$ws2.Cells.Item(48, 1).Value = "Cuba"
$ws2.Cells.Item(49, 1).Value = "Dinamarca"
$ws2.Cells.Item(50, 1).Value = "Dominica"
$ws2.Cells.Item(51, 1).Value = "Ecuador"
$ws2.Cells.Item(52, 1).Value = "Egipto"
$ws2.Cells.Item(53, 1).Value = "El Salvador"
$ws2.Cells.Item(54, 1).Value = "Emiratos Árabes Unidos"
$ws2.Cells.Item(55, 1).Value = "Eritrea"
$ws2.Cells.Item(56, 1).Value = "Eslovaquia"
$ws2.Cells.Item(57, 1).Value = "Eslovenia"
$ws2.Cells.Item(58, 1).Value = "España"
$ws2.Cells.Item(59, 1).Value = "Estados Unidos"
$ws2.Cells.Item(60, 1).Value = "Estonia"
$ws2.Cells.Item(61, 1).Value = "Etiopía"
$ws2.Cells.Item(62, 1).Value = "Filipinas"
$ws2.Cells.Item(63, 1).Value = "Finlandia"
$ws2.Cells.Item(64, 1).Value = "Fiyi"
$ws2.Cells.Item(65, 1).Value = "Francia"
$ws2.Cells.Item(66, 1).Value = "Gabón"
$ws2.Cells.Item(67, 1).Value = "Gambia"
$ws2.Cells.Item(68, 1).Value = "Georgia"
$ws2.Cells.Item(69, 1).Value = "Ghana"
$ws2.Cells.Item(70, 1).Value = "Granada"
$ws2.Cells.Item(71, 1).Value = "Grecia"
$ws2.Cells.Item(72, 1).Value = "Guatemala"
$ws2.Cells.Item(73, 1).Value = "Guinea"
$ws2.Cells.Item(74, 1).Value = "Guinea ecuatorial"
$ws2.Cells.Item(75, 1).Value = "Guinea-Bisáu"
$ws2.Cells.Item(76, 1).Value = "Haití"
$ws2.Cells.Item(77, 1).Value = "Honduras"
$ws2.Cells.Item(78, 1).Value = "Hungría"
$ws2.Cells.Item(79, 1).Value = "India"
$ws2.Cells.Item(80, 1).Value = "Indonesia"
$ws2.Cells.Item(81, 1).Value = "Irak"
$ws2.Cells.Item(82, 1).Value = "Irán"
$ws2.Cells.Item(83, 1).Value = "Irlanda"
$ws2.Cells.Item(84, 1).Value = "Islandia"
$ws2.Cells.Item(85, 1).Value = "Islas Marshall"
$ws2.Cells.Item(86, 1).Value = "Islas Salomón"
$ws2.Cells.Item(87, 1).Value = "Israel"
$ws2.Cells.Item(88, 1).Value = "Italia"
$ws2.Cells.Item(89, 1).Value = "Jamaica"
$ws2.Cells.Item(90, 1).Value = "Japón"
$ws2.Cells.Item(91, 1).Value = "Jordania"
$ws2.Cells.Item(92, 1).Value = "Kazajistán"
$ws2.Cells.Item(93, 1).Value = "Kenia"
$ws2.Cells.Item(94, 1).Value = "Kirguistán"
$ws2.Cells.Item(95, 1).Value = "Kiribati"
$ws2.Cells.Item(96, 1).Value = "Kuwait"
$ws2.Cells.Item(97, 1).Value = "Laos"
$ws2.Cells.Item(98, 1).Value = "Lesoto"
$ws2.Cells.Item(99, 1).Value = "Letonia"
$ws2.Cells.Item(100, 1).Value = "Líbano"
$ws2.Cells.Item(101, 1).Value = "Liberia"
$ws2.Cells.Item(102, 1).Value = "Libia"
$ws2.Cells.Item(103, 1).Value = "Liechtenstein"
$ws2.Cells.Item(104, 1).Value = "Lituania"
$ws2.Cells.Item(105, 1).Value = "Luxemburgo"
$ws2.Cells.Item(106, 1).Value = "Madagascar"
$ws2.Cells.Item(107, 1).Value = "Malasia"
$ws2.Cells.Item(108, 1).Value = "Malaui"
$ws2.Cells.Item(109, 1).Value = "Maldivas"
$ws2.Cells.Item(110, 1).Value = "Malí"
$ws2.Cells.Item(111, 1).Value = "Malta"
$ws2.Cells.Item(112, 1).Value = "Marruecos"
$ws2.Cells.Item(113, 1).Value = "Mauricio"
$ws2.Cells.Item(114, 1).Value = "Mauritania"
$ws2.Cells.Item(115, 1).Value = "México"
$ws2.Cells.Item(116, 1).Value = "Micronesia"
$ws2.Cells.Item(117, 1).Value = "Moldavia"
$ws2.Cells.Item(118, 1).Value = "Mónaco"
$ws2.Cells.Item(119, 1).Value = "Mongolia"
$ws2.Cells.Item(120, 1).Value = "Montenegro"
$ws2.Cells.Item(121, 1).Value = "Mozambique"
$ws2.Cells.Item(122, 1).Value = "Namibia"
$ws2.Cells.Item(123, 1).Value = "Nauru"
$ws2.Cells.Item(124, 1).Value = "Nepal"
$ws2.Cells.Item(125, 1).Value = "Nicaragua"
$ws2.Cells.Item(126, 1).Value = "Níger"
$ws2.Cells.Item(127, 1).Value = "Nigeria"
$ws2.Cells.Item(128, 1).Value = "Noruega"
$ws2.Cells.Item(129, 1).Value = "Nueva Zelanda"
$ws2.Cells.Item(130, 1).Value = "Omán"
$ws2.Cells.Item(131, 1).Value = "Países Bajos"
$ws2.Cells.Item(132, 1).Value = "Pakistán"
$ws2.Cells.Item(133, 1).Value = "Palaos"
$ws2.Cells.Item(134, 1).Value = "Panamá"
$ws2.Cells.Item(135, 1).Value = "Papúa Nueva Guinea"
$ws2.Cells.Item(136, 1).Value = "Paraguay"
$ws2.Cells.Item(137, 1).Value = "Perú"
$ws2.Cells.Item(138, 1).Value = "Polonia"
$ws2.Cells.Item(139, 1).Value = "Portugal"
$ws2.Cells.Item(140, 1).Value = "Reino Unido"
$ws2.Cells.Item(141, 1).Value = "República Centroafricana"
$ws2.Cells.Item(142, 1).Value = "República Checa"
$ws2.Cells.Item(143, 1).Value = "República de Macedonia"
$ws2.Cells.Item(144, 1).Value = "República del Congo"
$ws2.Cells.Item(145, 1).Value = "República Democrática del Congo"
$ws2.Cells.Item(146, 1).Value = "República Dominicana"
$ws2.Cells.Item(147, 1).Value = "República Sudafricana"
$ws2.Cells.Item(148, 1).Value = "Ruanda"
$ws2.Cells.Item(149, 1).Value = "Rumanía"
$ws2.Cells.Item(150, 1).Value = "Rusia"
$ws2.Cells.Item(151, 1).Value = "Samoa"
$ws2.Cells.Item(152, 1).Value = "San Cristóbal y Nieves"
$ws2.Cells.Item(153, 1).Value = "San Marino"
$ws2.Cells.Item(154, 1).Value = "San Vicente y las Granadinas"
$ws2.Cells.Item(155, 1).Value = "Santa Lucía"
$ws2.Cells.Item(156, 1).Value = "Santo Tomé y Príncipe"
$ws2.Cells.Item(157, 1).Value = "Senegal"
$ws2.Cells.Item(158, 1).Value = "Serbia"
$ws2.Cells.Item(159, 1).Value = "Seychelles"
$ws2.Cells.Item(160, 1).Value = "Sierra Leona"
$ws2.Cells.Item(161, 1).Value = "Singapur"
$ws2.Cells.Item(162, 1).Value = "Siria"
$ws2.Cells.Item(163, 1).Value = "Somalia"
$ws2.Cells.Item(164, 1).Value = "Sri Lanka"
$ws2.Cells.Item(165, 1).Value = "Suazilandia"
$ws2.Cells.Item(166, 1).Value = "Sudáfrica"
$ws2.Cells.Item(167, 1).Value = "Sudán"
$ws2.Cells.Item(168, 1).Value = "Sudán del Sur"
$ws2.Cells.Item(169, 1).Value = "Suecia"
$ws2.Cells.Item(170, 1).Value = "Suiza"
$ws2.Cells.Item(171, 1).Value = "Surinam"
$ws2.Cells.Item(172, 1).Value = "Tailandia"
$ws2.Cells.Item(173, 1).Value = "Tanzania"
$ws2.Cells.Item(174, 1).Value = "Tayikistán"
$ws2.Cells.Item(175, 1).Value = "Timor Oriental"
$ws2.Cells.Item(176, 1).Value = "Togo"
$ws2.Cells.Item(177, 1).Value = "Tonga"
$ws2.Cells.Item(178, 1).Value = "Trinidad y Tobago"
$ws2.Cells.Item(179, 1).Value = "Túnez"
$ws2.Cells.Item(180, 1).Value = "Turkmenistán"
$ws2.Cells.Item(181, 1).Value = "Turquía"
$ws2.Cells.Item(182, 1).Value = "Tuvalu"
$ws2.Cells.Item(183, 1).Value = "Ucrania"
$ws2.Cells.Item(184, 1).Value = "Uganda"
$ws2.Cells.Item(185, 1).Value = "Uruguay"
$ws2.Cells.Item(186, 1).Value = "Uzbekistán"
$ws2.Cells.Item(187, 1).Value = "Vanuatu"
$ws2.Cells.Item(188, 1).Value = "Venezuela"
$ws2.Cells.Item(189, 1).Value = "Vietnam"
$ws2.Cells.Item(190, 1).Value = "Yemen"
$ws2.Cells.Item(191, 1).Value = "Yibuti"
$ws2.Cells.Item(192, 1).Value = "Zambia"
$ws2.Cells.Item(193, 1).Value = "Zimbabue"

# States list (column B, rows 2-33)
$ws2.Cells.Item(2, 2).Value = "Aguascalientes"
$ws2.Cells.Item(3, 2).Value = "Baja California"
$ws2.Cells.Item(4, 2).Value = "Baja California Sur"
$ws2.Cells.Item(5, 2).Value = "Campeche"
$ws2.Cells.Item(6, 2).Value = "Chiapas"
$ws2.Cells.Item(7, 2).Value = "Chihuahua"
$ws2.Cells.Item(8, 2).Value = "Coahuila"
$ws2.Cells.Item(9, 2).Value = "Colima"
$ws2.Cells.Item(10, 2).Value = "Durango"
$ws2.Cells.Item(11, 2).Value = "Guanajuato"
$ws2.Cells.Item(12, 2).Value = "Guerrero"
$ws2.Cells.Item(13, 2).Value = "Hidalgo"
$ws2.Cells.Item(14, 2).Value = "Jalisco"
$ws2.Cells.Item(15, 2).Value = "Estado de México"
$ws2.Cells.Item(16, 2).Value = "Michoacán"
$ws2.Cells.Item(17, 2).Value = "Morelos"
$ws2.Cells.Item(18, 2).Value = "Nayarit"
$ws2.Cells.Item(19, 2).Value = "Nuevo León"
$ws2.Cells.Item(20, 2).Value = "Oaxaca"
$ws2.Cells.Item(21, 2).Value = "Puebla"
$ws2.Cells.Item(22, 2).Value = "Querétaro"
$ws2.Cells.Item(23, 2).Value = "Quintana Roo"
$ws2.Cells.Item(24, 2).Value = "San Luis Potosí"
$ws2.Cells.Item(25, 2).Value = "Sinaloa"
$ws2.Cells.Item(26, 2).Value = "Sonora"
$ws2.Cells.Item(27, 2).Value = "Tabasco"
$ws2.Cells.Item(28, 2).Value = "Tamaulipas"
$ws2.Cells.Item(29, 2).Value = "Tlaxcala"
$ws2.Cells.Item(30, 2).Value = "Veracruz"
$ws2.Cells.Item(31, 2).Value = "Yucatán"
$ws2.Cells.Item(32, 2).Value = "Zacatecas"
$ws2.Cells.Item(33, 2).Value = "Ciudad de México"

# Selections
$ws2.Range("R28").Select()
$ws1.Range("M24").Select()
